# Weekly fruit/vegetable price update:
# insert a new observation row right before the existing row 48 ("1a nueva(o)",
# Peru origin) and shift the rest of the "Zapallo"/"Camote" block down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 48 - everything from old row 48 down to old row 113
# (now 49..114) shifts down automatically, preserving formatting.
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with the new weekly record.
$ws.Range("A48").Value = 1
$ws.Range("B48").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C48").Value = "Arica y Parinacota"
$ws.Range("D48").Value = 45203
$ws.Range("E48").Value = 15
$ws.Range("F48").Value = 100112045
$ws.Range("G48").Value = "Zapallo"
$ws.Range("H48").Value = "Camote"
$ws.Range("I48").Value = "1a nueva(o)"
$ws.Range("J48").Value = 900
$ws.Range("K48").Value = 730
$ws.Range("L48").Value = 750
$ws.Range("M48").Value = 740
$ws.Range("N48").Value = "$/kilo (volumen en unidades)"
$ws.Range("O48").Value = "Perú"
$ws.Range("P48").Value = 740
$ws.Range("Q48").Value = 1
$ws.Range("R48").Value = "Hortaliza"
